# GB51251_2017_4_4_1 / GB50067 / GB50157 / GB50490 / GB50736 test-data update
# Adds:
#   - sheet1 (风口): new columns E "是否为室外风口" (bool, all FALSE) and
#     F "风量" (air volume) with per-row values.
#   - sheet4 (风机(不通过)) and sheet5 (风机(通过)): new column F "风机"
#     (air volume feeding the fan) with per-row values.
# Also restores the author's final navigation state (active sheet /
# selections) left over from editing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 风口 ("air outlet") -- add 是否为室外风口 (bool) + 风量 (number)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(1, 5).Value = "是否为室外风口"
$ws1.Cells.Item(1, 6).Value = "风量"

$outletFlow = @(
    @{Row=2;  F=1000},
    @{Row=3;  F=1000},
    @{Row=4;  F=1000},
    @{Row=5;  F=1000},
    @{Row=6;  F=1000},
    @{Row=7;  F=1500},
    @{Row=8;  F=1500},
    @{Row=9;  F=1000},
    @{Row=10; F=1000},
    @{Row=11; F=1000},
    @{Row=12; F=1000},
    @{Row=13; F=1000},
    @{Row=14; F=1000},
    @{Row=15; F=1000},
    @{Row=16; F=1000},
    @{Row=17; F=1000},
    @{Row=18; F=1000},
    @{Row=19; F=1000},
    @{Row=20; F=1000},
    @{Row=21; F=1500},
    @{Row=22; F=1500},
    @{Row=23; F=1000},
    @{Row=24; F=1000},
    @{Row=25; F=1000},
    @{Row=26; F=1500},
    @{Row=27; F=1000}
)

foreach ($item in $outletFlow) {
    $ws1.Cells.Item($item.Row, 5).Value = $false
    $ws1.Cells.Item($item.Row, 6).Value = $item.F
}

# ---------------------------------------------------------------------
# Sheet 4: 风机(不通过) -- add 风机 (number) column F
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(1, 6).Value = "风机"
$ws4.Cells.Item(2, 6).Value = 1000
$ws4.Cells.Item(3, 6).Value = 1000

# ---------------------------------------------------------------------
# Sheet 5: 风机(通过) -- add 风机 (number) column F
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Cells.Item(1, 6).Value = "风机"
$ws5.Cells.Item(2, 6).Value = 1000
$ws5.Cells.Item(3, 6).Value = 1000

# ---------------------------------------------------------------------
# Restore the author's final on-screen selections / active sheet.
# Sheet 1 ends up scrolled down with F28:F36 selected (just past the
# newly extended table), sheet 4 and sheet 5 both end with F1:F3
# selected, and sheet 5 (风机(通过)) is the sheet left active.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("F28:F36").Select()

$ws4.Activate()
$ws4.Range("F1:F3").Select()

$ws5.Activate()
$ws5.Range("F1:F3").Select()
